$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two new blank rows at 13:14. Old rows 13-18 shift to 15-20.
#    The new rows inherit row 12's formatting (B=s7, C/D/E=s4), which
#    already matches the desired format for the new row 13.
# ------------------------------------------------------------------
$ws.Rows("13:14").Insert()

# Row 14's B-column needs the "plain" percentage style (same as B11)
# instead of the inherited bordered one, so copy formats from B11.
$ws.Range("B11").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Append five new rows (21-25) at the bottom, duplicating row 20's
#    formatting (B=s8, C/D/E=s4) via copy + row-insert (this preserves
#    the exact style index without leaving stray values behind).
# ------------------------------------------------------------------
$ws.Rows("20").Copy()
$ws.Rows("21:25").Insert()
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Populate the "Dye" (A) column. Order matters: it controls the
#    order in which brand-new strings are appended to the shared
#    string table, which must match 10c,10d,12b,19a,22a,
#    Phenothiazine3,Phenothiazine4.
# ------------------------------------------------------------------
$ws.Range("A21").Value = "10c"
$ws.Range("A22").Value = "10d"
$ws.Range("A23").Value = "12b"
$ws.Range("A24").Value = "19a"
$ws.Range("A25").Value = "22a"

$ws.Range("A13").Value = "Phenothiazine3"
$ws.Range("A14").Value = "Phenothiazine4"

# ------------------------------------------------------------------
# 4. Fill in the experimental values for the new rows.
# ------------------------------------------------------------------
$ws.Range("B13").Value = 0.0243
$ws.Range("C13").Value = 0.559
$ws.Range("D13").Value = 5.76
$ws.Range("E13").Value = 0.71

$ws.Range("B14").Value = 0.0349
$ws.Range("C14").Value = 0.62
$ws.Range("D14").Value = 8.28
$ws.Range("E14").Value = 0.68

$ws.Range("B21").Value = 0.0541
$ws.Range("C21").Value = 0.605
$ws.Range("D21").Value = 12.9
$ws.Range("E21").Value = 0.7

$ws.Range("B22").Value = 0.050499999999999996
$ws.Range("C22").Value = 0.594
$ws.Range("D22").Value = 11.9
$ws.Range("E22").Value = 0.71

$ws.Range("B23").Value = 0.0231
$ws.Range("C23").Value = 0.505
$ws.Range("D23").Value = 6.7
$ws.Range("E23").Value = 0.68

$ws.Range("B24").Value = 0.0219
$ws.Range("C24").Value = 0.48
$ws.Range("D24").Value = 7.3
$ws.Range("E24").Value = 0.63

$ws.Range("B25").Value = 0.0441
$ws.Range("C25").Value = 0.566
$ws.Range("D25").Value = 11.1
$ws.Range("E25").Value = 0.7

# ------------------------------------------------------------------
# 5. Column A needs a custom width (closest reachable value to the
#    21.21875 target, given the engine's column-width quantization).
# ------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 20.3

# ------------------------------------------------------------------
# 6. Move the active selection to D15, matching the saved view state.
# ------------------------------------------------------------------
$ws.Range("D15").Select() | Out-Null
